$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149 - shifts existing rows 149:168 down to 150:169,
# preserving their data/formatting (matches the diff's row-shift pattern).
$ws.Rows("149:149").Insert()

# Populate the newly inserted row 149 with the new data record.
$ws.Range("A149").Value = 7
$ws.Range("B149").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C149").Value = "Ñuble"
$ws.Range("D149").Value = 44505
$ws.Range("E149").Value = 16
$ws.Range("F149").Value = 100112043
$ws.Range("G149").Value = "Pepino ensalada"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 120
$ws.Range("K149").Value = 9000
$ws.Range("L149").Value = 10000
$ws.Range("M149").Value = 9500
$ws.Range("N149").Value = "`$/caja 80 unidades"
$ws.Range("O149").Value = "Región del Maule"
$ws.Range("P149").Value = 119
$ws.Range("Q149").Value = 80
$ws.Range("R149").Value = "Hortaliza"
